$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" list (rows 16-22) is refreshed: the periods now run in
# ascending order (1911 -> 2005) instead of descending (2005 -> 1911), so
# each row's period label is mirrored around the middle row (19 / "2002",
# which stays put).
$ws.Range("E16").Value = "1911"
$ws.Range("E17").Value = "1912"
$ws.Range("E18").Value = "2001"
$ws.Range("E19").Value = "2002"
$ws.Range("E20").Value = "2003"
$ws.Range("E21").Value = "2004"
$ws.Range("E22").Value = "2005"

# The "Valor Mora" figures follow the same rows, so the value that used to
# belong to the 2005 period (now 1911, row 16) and the one that belonged to
# the 1911 period (now 2005, row 22) swap places.
$ws.Range("F16").Value = 60000
$ws.Range("F22").Value = 38000
